$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"). Copy H1's formatting first so the
# new header cells reuse the same bold/border/center header style, then set
# the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-25 for the new columns I ("I0") and J ("IF").
# I is 1 for every row except row 7 (2); J mirrors H except row 7 (H+1 = 4).
$iValues = @{
    2=1; 3=1; 4=1; 5=1; 6=1; 7=2; 8=1; 9=1; 10=1; 11=1; 12=1; 13=1; 14=1
    15=1; 16=1; 17=1; 18=1; 19=1; 20=1; 21=1; 22=1; 23=1; 24=1; 25=1
}
$jValues = @{
    2=5; 3=4; 4=7; 5=7; 6=4; 7=4; 8=5; 9=4; 10=6; 11=6; 12=3; 13=5; 14=6
    15=5; 16=5; 17=5; 18=4; 19=7; 20=5; 21=4; 22=3; 23=4; 24=3; 25=3
}

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
